{"js": "const replacements = [\n  [\"48\u00d721=\", \"23\u00d736=\"],\n  [\"56\u00d725=\", \"73\u00d791=\"],\n  [\"40\u00d770=\", \"37\u00d712=\"],\n  [\"95\u00d728=\", \"79\u00d711=\"],\n  [\"33\u00d753=\", \"87\u00d727=\"],\n  [\"89\u00d730=\", \"58\u00d789=\"],\n  [\"45\u00d728=\", \"80\u00d733=\"],\n  [\"81\u00d744=\", \"91\u00d711=\"],\n  [\"33\u00d719=\", \"96\u00d737=\"],\n  [\"85\u00d769=\", \"69\u00d746=\"],\n  [\"76\u00d794=\", \"59\u00d773=\"],\n  [\"70\u00d758=\", \"22\u00d754=\"],\n  [\"26\u00d748=\", \"12\u00d796=\"],\n  [\"23\u00d741=\", \"13\u00d768=\"],\n  [\"80\u00d766=\", \"62\u00d732=\"],\n  [\"84\u00d792=\", \"11\u00d742=\"],\n  [\"76\u00d722=\", \"71\u00d772=\"],\n  [\"35\u00d714=\", \"23\u00d767=\"],\n  [\"79\u00d762=\", \"73\u00d761=\"],\n  [\"62\u00d756=\", \"21\u00d726=\"],\n  [\"36\u00d791=\", \"39\u00d721=\"],\n  [\"83\u00d797=\", \"52\u00d751=\"],\n  [\"70\u00d772=\", \"67\u00d756=\"],\n  [\"79\u00d774=\", \"81\u00d736=\"],\n  [\"85\u00d746=\", \"43\u00d715=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"48\u00d721=\", \"23\u00d736=\"),\n    @(\"56\u00d725=\", \"73\u00d791=\"),\n    @(\"40\u00d770=\", \"37\u00d712=\"),\n    @(\"95\u00d728=\", \"79\u00d711=\"),\n    @(\"33\u00d753=\", \"87\u00d727=\"),\n    @(\"89\u00d730=\", \"58\u00d789=\"),\n    @(\"45\u00d728=\", \"80\u00d733=\"),\n    @(\"81\u00d744=\", \"91\u00d711=\"),\n    @(\"33\u00d719=\", \"96\u00d737=\"),\n    @(\"85\u00d769=\", \"69\u00d746=\"),\n    @(\"76\u00d794=\", \"59\u00d773=\"),\n    @(\"70\u00d758=\", \"22\u00d754=\"),\n    @(\"26\u00d748=\", \"12\u00d796=\"),\n    @(\"23\u00d741=\", \"13\u00d768=\"),\n    @(\"80\u00d766=\", \"62\u00d732=\"),\n    @(\"84\u00d792=\", \"11\u00d742=\"),\n    @(\"76\u00d722=\", \"71\u00d772=\"),\n    @(\"35\u00d714=\", \"23\u00d767=\"),\n    @(\"79\u00d762=\", \"73\u00d761=\"),\n    @(\"62\u00d756=\", \"21\u00d726=\"),\n    @(\"36\u00d791=\", \"39\u00d721=\"),\n    @(\"83\u00d797=\", \"52\u00d751=\"),\n    @(\"70\u00d772=\", \"67\u00d756=\"),\n    @(\"79\u00d774=\", \"81\u00d736=\"),\n    @(\"85\u00d746=\", \"43\u00d715=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
